# iron_native Word COM-interop script
# Applies: merges split/proofErr-wrapped runs back into single runs (no
# text change), and inserts a new paragraph about `npm run start`.

$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1) "Build in HTTPServer in python3" — collapse the spell-checked run split.
Replace-Exact "Build in HTTPServer in python3" "Build in HTTPServer in python3"

# 2) "First obtain the package using npm install . It will load the package.json "
Replace-Exact "First obtain the package using npm install . It will load the package.json " "First obtain the package using npm install . It will load the package.json "

# 3) "Inside /Client/App.js, you can config the HOST in line 7 for the backend url (Currently is pointing to localhost 8081)"
Replace-Exact "Inside /Client/App.js, you can config the HOST in line 7 for the backend url (Currently is pointing to localhost 8081)" "Inside /Client/App.js, you can config the HOST in line 7 for the backend url (Currently is pointing to localhost 8081)"

# 4) "Purchase orders that depends on USER and ITEM records"
Replace-Exact "Purchase orders that depends on USER and ITEM records" "Purchase orders that depends on USER and ITEM records"

# 5) "Please refer to WS_XXXXX.yml for more information"
Replace-Exact "Please refer to WS_XXXXX.yml for more information" "Please refer to WS_XXXXX.yml for more information"

# 6) "sing cryptography / cryptocode / simple-crypt will show invalid elf header in AWS Lambda (Under osx development). This project will use JWT for configuration encryption and decryption."
Replace-Exact "sing cryptography / cryptocode / simple-crypt will show invalid elf header in AWS Lambda (Under osx development). This project will use JWT for configuration encryption and decryption." "sing cryptography / cryptocode / simple-crypt will show invalid elf header in AWS Lambda (Under osx development). This project will use JWT for configuration encryption and decryption."

# 7) Insert a new paragraph "You can run local development via npm run start"
#    right after the blank paragraph that follows the "localhost 8081)"
#    paragraph, and before the next page-break paragraph.
$rng = $d.Content
$rng.Find.Execute("(Currently is pointing to localhost 8081)") | Out-Null
$para = $rng.Paragraphs(1)
$blankPara = $para.Next()
$blankPara.Range.InsertParagraphAfter()
$newPara = $blankPara.Next()
$newPara.Range.Text = "You can run local development via npm run start"
